$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new cells with shared-string values (B1="test1", C1="test2"),
# matching the author typing across the row from A1 (already "test").
$ws.Range("B1").Value = "test1"
$ws.Range("C1").Value = "test2"

# After typing into C1, Excel's selection moves on to the next cell (D1),
# which is what the saved sheetView/selection reflects.
$ws.Range("D1").Select()
